$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append two new data rows (121, 122) to the "Sheet 1" price history table,
# as produced by the authors' R script. Columns are:
#   A date | B volume | C high | D low | E open | F close | G adj_close | H ticker
# ---------------------------------------------------------------------------

# --- Row 121 ---------------------------------------------------------------
$ws.Cells.Item(121, 1).Value = 45496.2916666667
$ws.Cells.Item(121, 2).Value = 0
$ws.Cells.Item(121, 3).Value = 0.675000011920929
$ws.Cells.Item(121, 4).Value = 0.675000011920929
$ws.Cells.Item(121, 5).Value = 0.675000011920929
$ws.Cells.Item(121, 6).Value = 0.675000011920929
$ws.Cells.Item(121, 8).Value = "BWZ.MI"

# adj_close (col G) is stored as TEXT in the source workbook (matches the
# "close" value verbatim). Writing a numeric-looking string straight into
# Value gets auto-coerced back to a number, so build it as a text formula
# result and flatten it to a plain value via copy/paste-special - this
# keeps the cell's style untouched (no quote-prefix xf gets created).
$ws.Cells.Item(121, 7).Formula = '="0.675000011920929"'
$ws.Cells.Item(121, 7).Copy()
$ws.Cells.Item(121, 7).PasteSpecial(-4163)

# Match column A's date/time display format by copying A120's formatting.
$ws.Cells.Item(120, 1).Copy()
$ws.Cells.Item(121, 1).PasteSpecial(-4122)

# --- Row 122 ---------------------------------------------------------------
$ws.Cells.Item(122, 1).Value = 45497.6493634259
$ws.Cells.Item(122, 2).Value = 3020
$ws.Cells.Item(122, 3).Value = 0.680000007152557
$ws.Cells.Item(122, 4).Value = 0.680000007152557
$ws.Cells.Item(122, 5).Value = 0.680000007152557
$ws.Cells.Item(122, 6).Value = 0.680000007152557
$ws.Cells.Item(122, 8).Value = "BWZ.MI"

$ws.Cells.Item(122, 7).Formula = '="0.680000007152557"'
$ws.Cells.Item(122, 7).Copy()
$ws.Cells.Item(122, 7).PasteSpecial(-4163)

$ws.Cells.Item(120, 1).Copy()
$ws.Cells.Item(122, 1).PasteSpecial(-4122)

$excel.CutCopyMode = 0
